# Apply the commit's changes:
#  1. Rename the video file referenced in D2 from
#     "C:\Phase1\Videos\Men Will Be Men - edited.mp4" to
#     "C:\Phase1\Videos\demo.mp4"
#  2. Update the "Wait time (sec) after Video ends" value in E2 from
#     77 seconds to 93 seconds (stored as a fraction-of-a-day time value).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the video file path (D2).
$ws.Range("D2").Value = "C:\Phase1\Videos\demo.mp4"

# Update the wait time (E2) -- 93 seconds expressed as a fraction of a day,
# keeping the existing time number format on the cell.
$ws.Range("E2").Value = 93 / 86400
